$wb = $excel.ActiveWorkbook

# Map of row -> new "想去人数" (F column) value for the sheets that hold
# the conference listing data ("展览" and "全部类型").
$updates = @{
    6  = 1763
    7  = 31
    8  = 751
    9  = 367
    12 = 106
    16 = 121
    17 = 141
    18 = 4110
    20 = 26
    21 = 456
    22 = 392
    23 = 958
    28 = 1868
    29 = 54
    30 = 43
    31 = 80
    32 = 186
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }

    # Row 24 gets both the "想去人数" (F) and "最低票价" (G) columns updated.
    $ws.Cells.Item(24, 6).Value = 1188
    $ws.Cells.Item(24, 7).Value = 60
}
